$wb = $excel.ActiveWorkbook

$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A2").Value = "Version: $newVersion"
$wsAbout.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Tianan No. 1 Coal Mine, China, M1208, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

$wsData = $wb.Worksheets.Item("Boundaries and methane sources")
for ($r = 2; $r -le 7; $r++) {
    $wsData.Cells.Item($r, 19).Value = $newVersion
}
